# Add 4-departments title row (row 2) and supporting shared strings,
# then move the active selection to L5 (xls import function follow-up).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @("a", "b", "c", "d", "e", "f", "g", "h", "i", "end")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}

$ws.Range("L5").Select()
